$wb = $excel.ActiveWorkbook

# The "settings" sheet gains a new "version" column (control version support)
$ws = $wb.Worksheets.Item("settings")

$ws.Range("C1").Value = "version"
$ws.Range("C2").Value = 1

# Make "settings" the active sheet/tab, with C3 as the last selected cell
$ws.Activate()
$ws.Range("C3").Select()
